$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new "GHP Log" rows are being added to the drive_name/drive_path table:
#   - "SCR 3 - P85 Service" right after the existing "SCR 3 - PC 2" row
#   - "SCR 4 - P85 Service" right after the existing "SCR 4 - PC 2" row
#
# Insert both blank rows first (structural shift), then fill in the values.
$ws.Rows.Item(134).Insert()
$ws.Rows.Item(136).Insert()

# Fill the SCR 4 row (path first, then name).
$ws.Range("B136").Value = "\\10.214.86.219\d$\MES_Robust\GHP_19 4.0.12\P85SCR_L04S01\GHPService_P85SCR"
$ws.Range("A136").Value = "SCR 4 - P85 Service"

# Fill the SCR 3 row (name first, then path).
$ws.Range("A134").Value = "SCR 3 - P85 Service"
$ws.Range("B134").Value = "\\10.214.85.219\d$\MES_GHP Robust\GHP_18 4.0.3.0\P85SCR_L03S01\GHPService_P85SCR_L03S01"
